$wb = $excel.ActiveWorkbook

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 4682
$ws.Range("J58").Value = 8000
$ws.Range("L58").Value = 24000
$ws.Range("N58").Value = -24300

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7125.5
$ws.Range("I64").Value = 7001
$ws.Range("K64").Value = 7001
$ws.Range("M64").Value = -6753

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 7125.5
$ws.Range("I67").Value = 7001
$ws.Range("K67").Value = 7001
$ws.Range("M67").Value = -6143

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 55557396
$ws.Range("I86").Value = 76924930
$ws.Range("K86").Value = 76924930
$ws.Range("M86").Value = -76923807

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 55557396
$ws.Range("I89").Value = 76924930
$ws.Range("K89").Value = 384624650
$ws.Range("M89").Value = -384619034

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1774.3704
$ws.Range("I132").Value = 1751.3137
$ws.Range("J132").Value = 2166.3333
$ws.Range("K132").Value = 5253.9411
$ws.Range("L132").Value = 6498.999899999999
$ws.Range("M132").Value = -2723.9411
$ws.Range("N132").Value = -11558.9999

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2905.5908
$ws.Range("I137").Value = 3042.5
$ws.Range("K137").Value = 9127.5
$ws.Range("M137").Value = -6577.5

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2195965.8
$ws.Range("I32").Value = 2502921.2
$ws.Range("K32").Value = 2502921.2
$ws.Range("M32").Value = -2502634.2

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 42678.12
$ws.Range("I74").Value = 54260.95
$ws.Range("K74").Value = 54260.95
$ws.Range("M74").Value = -53386.95

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 42678.12
$ws.Range("I77").Value = 54260.95
$ws.Range("K77").Value = 271304.75
$ws.Range("M77").Value = -266936.75

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5209314
$ws.Range("I20").Value = 7247370.5
$ws.Range("K20").Value = 7247370.5
$ws.Range("M20").Value = -7247123.5

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 70886.47
$ws.Range("I86").Value = 169550.5
$ws.Range("K86").Value = 169550.5
$ws.Range("M86").Value = -168427.5

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 70886.47
$ws.Range("I89").Value = 169550.5
$ws.Range("K89").Value = 847752.5
$ws.Range("M89").Value = -842136.5

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1302.881
$ws.Range("J107").Value = 1788.3684
$ws.Range("L107").Value = 1788.3684
$ws.Range("N107").Value = -5628.3684

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2738.8333
$ws.Range("I122").Value = 2518.6875
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 7556.0625
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -5106.0625
$ws.Range("N122").Value = -18400

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7160.6924
$ws.Range("I132").Value = 4407.5386
$ws.Range("K132").Value = 13222.6158
$ws.Range("M132").Value = -10692.6158

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3066.5557
$ws.Range("I113").Value = 678.75
$ws.Range("J113").Value = 4071.9473
$ws.Range("K113").Value = 2036.25
$ws.Range("L113").Value = 12215.8419
$ws.Range("M113").Value = 133.75
$ws.Range("N113").Value = -16555.8419

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 11977097
$ws.Range("I129").Value = 960.75
$ws.Range("J129").Value = 27945278
$ws.Range("K129").Value = 2882.25
$ws.Range("L129").Value = 83835834
$ws.Range("M129").Value = 2117.75
$ws.Range("N129").Value = -83845834

# GSM row 62
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 110000
$ws.Range("J62").Value = 110000
$ws.Range("L62").Value = 110000
$ws.Range("N62").Value = -111372

# GSM row 65
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 110000
$ws.Range("J65").Value = 110000
$ws.Range("L65").Value = 330000
$ws.Range("N65").Value = -336864

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9348.6
$ws.Range("I70").Value = 6990
$ws.Range("J70").Value = 9938.25
$ws.Range("K70").Value = 6990
$ws.Range("L70").Value = 9938.25
$ws.Range("M70").Value = -6720
$ws.Range("N70").Value = -10478.25

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 9348.6
$ws.Range("I73").Value = 6990
$ws.Range("J73").Value = 9938.25
$ws.Range("K73").Value = 6990
$ws.Range("L73").Value = 9938.25
$ws.Range("M73").Value = -6054
$ws.Range("N73").Value = -11810.25

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3659.6924
$ws.Range("I80").Value = 3354.9092
$ws.Range("K80").Value = 3354.9092
$ws.Range("M80").Value = -2356.9092

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3659.6924
$ws.Range("I83").Value = 3354.9092
$ws.Range("K83").Value = 16774.546
$ws.Range("M83").Value = -11782.546

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4226.1304
$ws.Range("I7").Value = 3013.1333
$ws.Range("K7").Value = 3013.1333
$ws.Range("M7").Value = -2901.1333

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4185.4062
$ws.Range("I22").Value = 1969.3684
$ws.Range("J22").Value = 7424.231
$ws.Range("K22").Value = 1969.3684
$ws.Range("L22").Value = 7424.231
$ws.Range("M22").Value = -1674.3684
$ws.Range("N22").Value = -8014.231

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 4185.4062
$ws.Range("I27").Value = 1969.3684
$ws.Range("J27").Value = 7424.231
$ws.Range("K27").Value = 1969.3684
$ws.Range("L27").Value = 7424.231
$ws.Range("M27").Value = -1862.3684
$ws.Range("N27").Value = -7638.231

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1551.48
$ws.Range("I46").Value = 848
$ws.Range("K46").Value = 848
$ws.Range("M46").Value = -660

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 529.5454999999999
$ws.Range("I55").Value = 113
$ws.Range("J55").Value = 685.75
$ws.Range("K55").Value = 113
$ws.Range("L55").Value = 685.75
$ws.Range("M55").Value = 60
$ws.Range("N55").Value = -1031.75

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3276.2942
$ws.Range("I93").Value = 2791.2307
$ws.Range("K93").Value = 2791.2307
$ws.Range("M93").Value = -1543.2307

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4226.1304
$ws.Range("I126").Value = 3013.1333
$ws.Range("K126").Value = 9039.3999
$ws.Range("M126").Value = -6569.3999

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6360.357
$ws.Range("I132").Value = 3361.7693
$ws.Range("J132").Value = 11233.0625
$ws.Range("K132").Value = 10085.3079
$ws.Range("L132").Value = 33699.1875
$ws.Range("M132").Value = -7555.3079
$ws.Range("N132").Value = -38759.1875

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9663.677
$ws.Range("I136").Value = 4448.5
$ws.Range("K136").Value = 13345.5
$ws.Range("M136").Value = -10795.5

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 13623886
$ws.Range("I122").Value = 18668134
$ws.Range("K122").Value = 56004402
$ws.Range("M122").Value = -56001952

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 37041836
$ws.Range("I126").Value = 62503932
$ws.Range("J126").Value = 6055.364
$ws.Range("K126").Value = 187511796
$ws.Range("L126").Value = 18166.092
$ws.Range("M126").Value = -187509326
$ws.Range("N126").Value = -23106.092

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 33957.25
$ws.Range("I132").Value = 3094
$ws.Range("K132").Value = 9282
$ws.Range("M132").Value = -6752

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 38530.43
$ws.Range("I136").Value = 2018.5454
$ws.Range("K136").Value = 6055.6362
$ws.Range("M136").Value = -3505.6362
